# Update degree abbreviations to include periods, per commit message:
# "Updated author pages to have '.' in the degree abbreviations"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Ph.D. Forestry"
$ws.Range("A17").Value = "M.Sc. Biology"
$ws.Range("A18").Value = "B.Sc. Environmental Science"

$ws.Range("A18").Select()
